$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.251.31"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.775.31"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "625.39"
$ws.Range("E5").Value = "  +4.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.93"
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.774.49"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.460"
$ws.Range("E11").Value = "  +3.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.73"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.83"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.414.23"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.774.96"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.265.97"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.67"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.08"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.24"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.58"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("E24").Value = "  +3.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.15"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("E27").Value = "  +4.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.927.51"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.25"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.18"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.82"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.166"
$ws.Range("E36").Value = "  +14.25%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.729.48"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.98"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  +7.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.82"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.968"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.297"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.93"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.01"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.74"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.91"
$ws.Range("E49").Value = "  +4.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.41"
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("E51").Value = "  -0.17%  "
